$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    paragraph (the Heading1 "Play Caribbean Voyage Free Slot Game -
#    Funta Gaming" title), and before the "Gameplay" heading.
# ---------------------------------------------------------------------
$metaXml = '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
           '<w:r><w:t>: Set sail for adventure in Caribbean Voyage, the exciting pirate-themed slot from Funta Gaming. Play for free or real money and win big with unique bonus features.</w:t></w:r></w:p>' +
           '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Gameplay</w:t></w:r></w:p>' +
           '</w:body></w:document>'

# Replacing the "Gameplay" heading paragraph's content with itself plus a
# new paragraph in front of it effectively inserts that new paragraph
# right before "Gameplay" (i.e. right after the title paragraph).
$gameplayPara = $d.Paragraphs(2)
[void]$gameplayPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) Remove the trailing "Play Caribbean Voyage Free Slot Game - Funta
#    Gaming" paragraph (the bold one near the end of the document) and
#    replace the text of the following italic paragraph with the new
#    image-prompt copy.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$playAgainPara = $d.Paragraphs($count - 1)
$playAgainPara.Range.Delete()

$imagePromptPara = $d.Paragraphs($d.Paragraphs.Count)
$imgRange = $d.Range($imagePromptPara.Range.Start, $imagePromptPara.Range.End)
$imgRange.Text = 'Create a fun and vibrant feature image for Funta Gaming''s latest slot game, Caribbean Voyage. The image should be in cartoon style and feature a happy Maya warrior with glasses, celebrating his winnings from the game. The warrior should be surrounded by bright and bold Caribbean inspired graphics such as palm trees, treasure chests, and sea creatures. The image should also include text that reads "Caribbean Voyage: The Ultimate Pirate Adventure!" in a fun and playful font. The overall vibe of the image should be exciting and adventurous, conveying the thrill of the game.'
